$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Col = 2;  Old = "收入（千美元）";             New = "收入 (`$K)" },
    @{ Col = 3;  Old = "销货成本（千美元）";         New = "所售货物成本 (`$K)" },
    @{ Col = 4;  Old = "毛利率 (%)";                New = "毛利润率 (%)" },
    @{ Col = 5;  Old = "运营费用（千美元）";         New = "运营费用 (`$K)" },
    @{ Col = 6;  Old = "息税折旧摊销前利润（千美元）"; New = "EBITDA  (`$K)" },
    @{ Col = 7;  Old = "利息支出（千美元）";         New = "利息支出 (`$K)" },
    @{ Col = 8;  Old = "税前收益（千美元）";         New = "税前收益 (`$K)" },
    @{ Col = 9;  Old = "净收入（千美元）";           New = "净收入 (`$K)" },
    @{ Col = 10; Old = "总资产（千美元）";           New = "总资产 (`$K)" },
    @{ Col = 11; Old = "总负债（千美元）";           New = "总负债 (`$K)" },
    @{ Col = 12; Old = "股东权益（千美元）";         New = "股东权益 (`$K)" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell(1, $item.Col)
    $r = $cell.Range
    $find = $r.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Bold = 1
    $find.Execute($item.Old, $true, $false, $false, $false, $false, $true, 0, $false, $item.New, 1) | Out-Null
}
